$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 144
$ws.Cells.Item(144,1).Value = 142
$ws.Cells.Item(144,2).Value = 7913218
$ws.Cells.Item(144,3).Value = 'Uruguay Primera División'
$ws.Cells.Item(144,4).Value = 'Uruguay Apertura'
$ws.Cells.Item(144,5).Value = 45359.85416666666
$ws.Cells.Item(144,6).Value = 'Miramar Misiones'
$ws.Cells.Item(144,7).Value = 'Atletico Fenix Montevideo'
$ws.Cells.Item(144,8).Value = 1
$ws.Cells.Item(144,9).Value = 1
$ws.Cells.Item(144,10).Value = 'D'
$ws.Cells.Item(144,11).Value = 2.625
$ws.Cells.Item(144,12).Value = 3
$ws.Cells.Item(144,13).Value = 2.6
$ws.Cells.Item(144,14).Value = 2.375
$ws.Cells.Item(144,15).Value = 3
$ws.Cells.Item(144,16).Value = 2.875
$ws.Cells.Item(144,17).Value = -0.25
$ws.Cells.Item(144,18).Value = 2.1
$ws.Cells.Item(144,19).Value = 1.775
$ws.Cells.Item(144,20).Value = 2.25
$ws.Cells.Item(144,21).Value = 2
$ws.Cells.Item(144,22).Value = 1.85
$ws.Cells.Item(144,23).Value = -1
$ws.Cells.Item(144,24).Value = 2
$ws.Cells.Item(144,25).Value = -1
$ws.Cells.Item(144,26).Value = -0.5
$ws.Cells.Item(144,27).Value = 0.3875
$ws.Cells.Item(144,28).Value = -0.5
$ws.Cells.Item(144,29).Value = 0.425

# Row 145
$ws.Cells.Item(145,1).Value = 143
$ws.Cells.Item(145,2).Value = 7913219
$ws.Cells.Item(145,3).Value = 'Uruguay Primera División'
$ws.Cells.Item(145,4).Value = 'Uruguay Apertura'
$ws.Cells.Item(145,5).Value = 45360.6875
$ws.Cells.Item(145,6).Value = 'CA River Plate'
$ws.Cells.Item(145,7).Value = 'Boston River'
$ws.Cells.Item(145,11).Value = 2.375
$ws.Cells.Item(145,12).Value = 3.1
$ws.Cells.Item(145,13).Value = 2.9
$ws.Cells.Item(145,14).Value = 2.45
$ws.Cells.Item(145,15).Value = 3
$ws.Cells.Item(145,16).Value = 2.875
$ws.Cells.Item(145,17).Value = 0
$ws.Cells.Item(145,18).Value = 1.775
$ws.Cells.Item(145,19).Value = 2.1
$ws.Cells.Item(145,20).Value = 2
$ws.Cells.Item(145,21).Value = 1.8
$ws.Cells.Item(145,22).Value = 2.05
$ws.Cells.Item(145,23).Value = 0
$ws.Cells.Item(145,24).Value = 0
$ws.Cells.Item(145,25).Value = 0
$ws.Cells.Item(145,26).Value = 0
$ws.Cells.Item(145,27).Value = 0

# Row 146
$ws.Cells.Item(146,1).Value = 144
$ws.Cells.Item(146,2).Value = 7913223
$ws.Cells.Item(146,3).Value = 'Uruguay Primera División'
$ws.Cells.Item(146,4).Value = 'Uruguay Apertura'
$ws.Cells.Item(146,5).Value = 45360.8125
$ws.Cells.Item(146,6).Value = 'Penarol'
$ws.Cells.Item(146,7).Value = 'Cerro'
$ws.Cells.Item(146,11).Value = 1.5
$ws.Cells.Item(146,12).Value = 4
$ws.Cells.Item(146,13).Value = 6
$ws.Cells.Item(146,14).Value = 1.444
$ws.Cells.Item(146,15).Value = 4.2
$ws.Cells.Item(146,16).Value = 7
$ws.Cells.Item(146,17).Value = -1
$ws.Cells.Item(146,18).Value = 1.775
$ws.Cells.Item(146,19).Value = 2.1
$ws.Cells.Item(146,20).Value = 2.25
$ws.Cells.Item(146,21).Value = 1.925
$ws.Cells.Item(146,22).Value = 1.925
$ws.Cells.Item(146,23).Value = 0
$ws.Cells.Item(146,24).Value = 0
$ws.Cells.Item(146,25).Value = 0
$ws.Cells.Item(146,26).Value = 0
$ws.Cells.Item(146,27).Value = 0

# Row 147
$ws.Cells.Item(147,1).Value = 145
$ws.Cells.Item(147,2).Value = 7913221
$ws.Cells.Item(147,3).Value = 'Uruguay Primera División'
$ws.Cells.Item(147,4).Value = 'Uruguay Apertura'
$ws.Cells.Item(147,5).Value = 45361.41666666666
$ws.Cells.Item(147,6).Value = 'Liverpool Montevideo'
$ws.Cells.Item(147,7).Value = 'Danubio'
$ws.Cells.Item(147,11).Value = 2.2
$ws.Cells.Item(147,12).Value = 3.2
$ws.Cells.Item(147,13).Value = 3.2
$ws.Cells.Item(147,14).Value = 2.2
$ws.Cells.Item(147,15).Value = 3.2
$ws.Cells.Item(147,16).Value = 3.25
$ws.Cells.Item(147,17).Value = -0.25
$ws.Cells.Item(147,18).Value = 1.925
$ws.Cells.Item(147,19).Value = 1.925
$ws.Cells.Item(147,20).Value = 2.25
$ws.Cells.Item(147,21).Value = 1.975
$ws.Cells.Item(147,22).Value = 1.875
$ws.Cells.Item(147,23).Value = 0
$ws.Cells.Item(147,24).Value = 0
$ws.Cells.Item(147,25).Value = 0
$ws.Cells.Item(147,26).Value = 0
$ws.Cells.Item(147,27).Value = 0

# Row 148
$ws.Cells.Item(148,1).Value = 146
$ws.Cells.Item(148,2).Value = 7913222
$ws.Cells.Item(148,3).Value = 'Uruguay Primera División'
$ws.Cells.Item(148,4).Value = 'Uruguay Apertura'
$ws.Cells.Item(148,5).Value = 45361.6875
$ws.Cells.Item(148,6).Value = 'Club Atletico Progreso'
$ws.Cells.Item(148,7).Value = 'Racing Club de Montevideo'
$ws.Cells.Item(148,11).Value = 2.75
$ws.Cells.Item(148,12).Value = 3
$ws.Cells.Item(148,13).Value = 2.625
$ws.Cells.Item(148,14).Value = 2.75
$ws.Cells.Item(148,15).Value = 3
$ws.Cells.Item(148,16).Value = 2.625
$ws.Cells.Item(148,17).Value = 0
$ws.Cells.Item(148,18).Value = 2
$ws.Cells.Item(148,19).Value = 1.85
$ws.Cells.Item(148,20).Value = 2.25
$ws.Cells.Item(148,21).Value = 2.025
$ws.Cells.Item(148,22).Value = 1.825
$ws.Cells.Item(148,23).Value = 0
$ws.Cells.Item(148,24).Value = 0
$ws.Cells.Item(148,25).Value = 0
$ws.Cells.Item(148,26).Value = 0
$ws.Cells.Item(148,27).Value = 0

# Row 149
$ws.Cells.Item(149,1).Value = 147
$ws.Cells.Item(149,2).Value = 7913220
$ws.Cells.Item(149,3).Value = 'Uruguay Primera División'
$ws.Cells.Item(149,4).Value = 'Uruguay Apertura'
$ws.Cells.Item(149,5).Value = 45361.8125
$ws.Cells.Item(149,6).Value = 'Defensor Sporting'
$ws.Cells.Item(149,7).Value = 'Nacional De Football'
$ws.Cells.Item(149,11).Value = 2.375
$ws.Cells.Item(149,12).Value = 3.2
$ws.Cells.Item(149,13).Value = 2.875
$ws.Cells.Item(149,14).Value = 2.2
$ws.Cells.Item(149,15).Value = 3.2
$ws.Cells.Item(149,16).Value = 3.2
$ws.Cells.Item(149,17).Value = -0.25
$ws.Cells.Item(149,18).Value = 1.975
$ws.Cells.Item(149,19).Value = 1.875
$ws.Cells.Item(149,20).Value = 2.25
$ws.Cells.Item(149,21).Value = 1.8
$ws.Cells.Item(149,22).Value = 2.05
$ws.Cells.Item(149,23).Value = 0
$ws.Cells.Item(149,24).Value = 0
$ws.Cells.Item(149,25).Value = 0
$ws.Cells.Item(149,26).Value = 0
$ws.Cells.Item(149,27).Value = 0

# Row 150
$ws.Cells.Item(150,1).Value = 148
$ws.Cells.Item(150,2).Value = 7913224
$ws.Cells.Item(150,3).Value = 'Uruguay Primera División'
$ws.Cells.Item(150,4).Value = 'Uruguay Apertura'
$ws.Cells.Item(150,5).Value = 45362.6875
$ws.Cells.Item(150,6).Value = 'Cerro Largo'
$ws.Cells.Item(150,7).Value = 'Rampla Juniors'
$ws.Cells.Item(150,11).Value = 1.909
$ws.Cells.Item(150,12).Value = 3.2
$ws.Cells.Item(150,13).Value = 4
$ws.Cells.Item(150,14).Value = 1.8
$ws.Cells.Item(150,15).Value = 3.3
$ws.Cells.Item(150,16).Value = 4.5
$ws.Cells.Item(150,17).Value = -0.5
$ws.Cells.Item(150,18).Value = 1.825
$ws.Cells.Item(150,19).Value = 2.025
$ws.Cells.Item(150,20).Value = 2.25
$ws.Cells.Item(150,21).Value = 2
$ws.Cells.Item(150,22).Value = 1.85
$ws.Cells.Item(150,23).Value = 0
$ws.Cells.Item(150,24).Value = 0
$ws.Cells.Item(150,25).Value = 0
$ws.Cells.Item(150,26).Value = 0
$ws.Cells.Item(150,27).Value = 0

# Row 151
$ws.Cells.Item(151,1).Value = 149
$ws.Cells.Item(151,2).Value = 7913225
$ws.Cells.Item(151,3).Value = 'Uruguay Primera División'
$ws.Cells.Item(151,4).Value = 'Uruguay Apertura'
$ws.Cells.Item(151,5).Value = 45362.79166666666
$ws.Cells.Item(151,6).Value = 'Montevideo Wanderers'
$ws.Cells.Item(151,7).Value = 'Deportivo Maldonado'
$ws.Cells.Item(151,11).Value = 2.375
$ws.Cells.Item(151,12).Value = 3.3
$ws.Cells.Item(151,13).Value = 2.75
$ws.Cells.Item(151,14).Value = 2.6
$ws.Cells.Item(151,15).Value = 3.3
$ws.Cells.Item(151,16).Value = 2.5
$ws.Cells.Item(151,17).Value = 0
$ws.Cells.Item(151,18).Value = 1.975
$ws.Cells.Item(151,19).Value = 1.875
$ws.Cells.Item(151,20).Value = 2.25
$ws.Cells.Item(151,21).Value = 1.95
$ws.Cells.Item(151,22).Value = 1.9
$ws.Cells.Item(151,23).Value = 0
$ws.Cells.Item(151,24).Value = 0
$ws.Cells.Item(151,25).Value = 0
$ws.Cells.Item(151,26).Value = 0
$ws.Cells.Item(151,27).Value = 0

# Fix formatting for the two brand-new rows (150, 151): copy id/date styles
$ws.Range("A149").Copy()
$ws.Range("A150:A151").PasteSpecial(-4122)
$ws.Range("E149").Copy()
$ws.Range("E150:E151").PasteSpecial(-4122)
$excel.CutCopyMode = $false
